$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cases")
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$axis = $chart.Axes(2)
$axis.DisplayUnit = 9
Write-Host "DisplayUnit: $($axis.DisplayUnit)"
